# Auto-generated Excel COM-interop script
# Applies market-data price/profit refresh updates across ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 5020
$ws.Range("I6").Value = 5020
$ws.Range("K6").Value = 15060
$ws.Range("M6").Value = -14948

$ws.Range("H15").Value = 17244548
$ws.Range("I15").Value = 17244548
$ws.Range("K15").Value = 51733644
$ws.Range("M15").Value = -51733475

$ws.Range("H74").Value = 5194962.5
$ws.Range("J74").Value = 5250
$ws.Range("L74").Value = 5250
$ws.Range("N74").Value = -7122

$ws.Range("H76").Value = 183339000
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H77").Value = 5194962.5
$ws.Range("J77").Value = 5250
$ws.Range("L77").Value = 26250
$ws.Range("N77").Value = -35610

$ws.Range("H79").Value = 183339000
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H106").Value = 3632.7778
$ws.Range("I106").Value = 3385
$ws.Range("J106").Value = 4500
$ws.Range("K106").Value = 3385
$ws.Range("L106").Value = 4500
$ws.Range("M106").Value = -2754
$ws.Range("N106").Value = -5762

$ws.Range("H113").Value = 2782.7727
$ws.Range("I113").Value = 2671.5334
$ws.Range("J113").Value = 3021.1428
$ws.Range("K113").Value = 2671.5334
$ws.Range("L113").Value = 3021.1428
$ws.Range("M113").Value = 582.4666000000002
$ws.Range("N113").Value = -9529.1428

$ws.Range("H125").Value = 4493.909
$ws.Range("I125").Value = 1443.3
$ws.Range("J125").Value = 35000
$ws.Range("K125").Value = 12989.7
$ws.Range("L125").Value = 315000
$ws.Range("M125").Value = -10529.7
$ws.Range("N125").Value = -319920

$ws.Range("H137").Value = 786.3823
$ws.Range("I137").Value = 726.6786
$ws.Range("J137").Value = 1065
$ws.Range("K137").Value = 2180.0358
$ws.Range("L137").Value = 3195
$ws.Range("M137").Value = 369.9642000000003
$ws.Range("N137").Value = -8295

$ws.Range("H138").Value = 2776.9192
$ws.Range("I138").Value = 734.4
$ws.Range("J138").Value = 3466.9595
$ws.Range("K138").Value = 2203.2
$ws.Range("L138").Value = 10400.8785
$ws.Range("M138").Value = 2936.8
$ws.Range("N138").Value = -20680.8785

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4804.8506
$ws.Range("I32").Value = 4099.283
$ws.Range("K32").Value = 4099.283
$ws.Range("M32").Value = -3812.283

$ws.Range("H74").Value = 696.413
$ws.Range("I74").Value = 618.59375
$ws.Range("J74").Value = 874.2857
$ws.Range("K74").Value = 618.59375
$ws.Range("L74").Value = 874.2857
$ws.Range("M74").Value = 255.40625
$ws.Range("N74").Value = -2622.2857

$ws.Range("H77").Value = 696.413
$ws.Range("I77").Value = 618.59375
$ws.Range("J77").Value = 874.2857
$ws.Range("K77").Value = 3092.96875
$ws.Range("L77").Value = 4371.4285
$ws.Range("M77").Value = 1275.03125
$ws.Range("N77").Value = -13107.4285

$ws.Range("H132").Value = 1960.1892
$ws.Range("I132").Value = 1457.1428
$ws.Range("J132").Value = 3525.2222
$ws.Range("K132").Value = 4371.428400000001
$ws.Range("L132").Value = 10575.6666
$ws.Range("M132").Value = -1841.428400000001
$ws.Range("N132").Value = -15635.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4432.6787
$ws.Range("I105").Value = 4471.8687
$ws.Range("J105").Value = 4349.9443
$ws.Range("K105").Value = 4471.8687
$ws.Range("L105").Value = 4349.9443
$ws.Range("M105").Value = -2724.8687
$ws.Range("N105").Value = -7843.9443

$ws.Range("H134").Value = 21211.54
$ws.Range("I134").Value = 1160.0698
$ws.Range("J134").Value = 144384.86
$ws.Range("K134").Value = 3480.2094
$ws.Range("L134").Value = 433154.58
$ws.Range("M134").Value = -945.2093999999997
$ws.Range("N134").Value = -438224.58

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3504.7693
$ws.Range("I58").Value = 621.2593000000001
$ws.Range("J58").Value = 9992.666999999999
$ws.Range("K58").Value = 621.2593000000001
$ws.Range("L58").Value = 9992.666999999999
$ws.Range("M58").Value = -418.2593000000001
$ws.Range("N58").Value = -10398.667

$ws.Range("H86").Value = 3374.1667
$ws.Range("I86").Value = 3212.8572
$ws.Range("K86").Value = 3212.8572
$ws.Range("M86").Value = -2089.8572

$ws.Range("H89").Value = 3374.1667
$ws.Range("I89").Value = 3212.8572
$ws.Range("K89").Value = 16064.286
$ws.Range("M89").Value = -10448.286

$ws.Range("H99").Value = 2686.1724
$ws.Range("I99").Value = 2435.2942
$ws.Range("J99").Value = 3041.5833
$ws.Range("K99").Value = 2435.2942
$ws.Range("L99").Value = 3041.5833
$ws.Range("M99").Value = -937.2941999999998
$ws.Range("N99").Value = -6037.5833

$ws.Range("H126").Value = 2686.1724
$ws.Range("I126").Value = 2435.2942
$ws.Range("J126").Value = 3041.5833
$ws.Range("K126").Value = 7305.882599999999
$ws.Range("L126").Value = 9124.749899999999
$ws.Range("M126").Value = -4835.882599999999
$ws.Range("N126").Value = -14064.7499

$ws.Range("H136").Value = 3504.7693
$ws.Range("I136").Value = 621.2593000000001
$ws.Range("J136").Value = 9992.666999999999
$ws.Range("K136").Value = 1863.7779
$ws.Range("L136").Value = 29978.001
$ws.Range("M136").Value = 686.2221
$ws.Range("N136").Value = -35078.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 8000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 8000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -8540

$ws.Range("H73").Value = 8000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 8000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 8000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -9872

$ws.Range("H126").Value = 2166.5789
$ws.Range("I126").Value = 2034.4375
$ws.Range("J126").Value = 2871.3333
$ws.Range("K126").Value = 6103.3125
$ws.Range("L126").Value = 8613.999899999999
$ws.Range("M126").Value = -3633.3125
$ws.Range("N126").Value = -13553.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 10680.8
$ws.Range("I122").Value = 12576
$ws.Range("J122").Value = 3100
$ws.Range("K122").Value = 37728
$ws.Range("L122").Value = 9300
$ws.Range("M122").Value = -35278
$ws.Range("N122").Value = -14200

$ws.Range("H132").Value = 2484.4897
$ws.Range("I132").Value = 2348.6191
$ws.Range("J132").Value = 3299.7144
$ws.Range("K132").Value = 7045.8573
$ws.Range("L132").Value = 9899.143199999999
$ws.Range("M132").Value = -4515.8573
$ws.Range("N132").Value = -14959.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2012.65
$ws.Range("I136").Value = 1841.2413
$ws.Range("J136").Value = 2464.5454
$ws.Range("K136").Value = 5523.7239
$ws.Range("L136").Value = 7393.6362
$ws.Range("M136").Value = -2973.7239
$ws.Range("N136").Value = -12493.6362
